# Zeitaufteilung.xlsx - Zeiten aktualisiert und in allen Kapiteln eingetragen, Gendern
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- 1. Update the time-allocation figures that changed ---
$ws.Range("C4").Value = 15
$ws.Range("B5").Value = 90
$ws.Range("B6").Value = 75
$ws.Range("B9").Value = 60
$ws.Range("B10").Value = 0

# --- 2. Row 8 previously carried an explicit row-level custom format; ---
#        drop it (while keeping the taller row height) before re-applying ---
#        the same top-alignment the rest of the sheet now uses. ---
$rowHeight8 = $ws.Rows("8:8").RowHeight
$ws.Rows("8:8").ClearFormats()
$ws.Rows("8:8").RowHeight = $rowHeight8

# --- 3. Align the content of every populated cell to the top (the whole ---
#        table was reformatted this way) ---
$cells = @("A1","B1","C1","D1","E1","F1","A2","B2","D2","A3","B3","D3","A4","B4","C4","D4","E4","A5","B5","C5","D5","E5","F5","A6","B6","D6","A7","B7","C7","D7","E7","A8","B8","C8","D8","E8","F8","A9","B9","D9","A10","B10","D10","B12","C12","D12","A14","B14","C14","E14")
foreach ($addr in $cells) {
    $ws.Range($addr).VerticalAlignment = -4160
}
# E8 keeps its wrapped text on top of the top alignment
$ws.Range("E8").WrapText = $true

# --- 4. Column D grew slightly wider to fit the new values ---
$ws.Columns("D").ColumnWidth = 6.6667

# --- 5. Update the active selection ---
$ws.Range("B7").Select()
